$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.685.63'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +2.10%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.256.40'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.99%  '
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.82'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.43%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.632'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.38%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '76.90'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +3.24%  '
$ws.Range("E8").Value = '  +0.10%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.630'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.43%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '45.38'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +14.31%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0955'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.95%  '
$ws.Range("E12").Value = '  +2.81%  '
$ws.Range("E13").Value = '  -0.15%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '14.76'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.18%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.865'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.93%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.257.71'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +1.33%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '42.510.09'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +1.86%  '
$ws.Range("E18").Value = '  +4.06%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.22'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +1.83%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '72.24'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +1.32%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.82'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +51.22%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.27'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +1.19%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '232.71'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +1.16%  '
$ws.Range("E24").Value = '  +6.26%  '
$ws.Range("E25").Value = '  +0.08%  '
$ws.Range("E26").Value = '  -2.19%  '
$ws.Range("E27").Value = '  +1.16%  '
$ws.Range("E28").Value = '  +5.53%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '167.39'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.79%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '20.72'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +1.20%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0830'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.94%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '32.33'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -6.12%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.61'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +15.25%  '
$ws.Range("E34").Value = '  +0.89%  '
$ws.Range("E35").Value = '  +0.05%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.68'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +1.62%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0317'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +6.40%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '14.50'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +7.26%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.20'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.07%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.83'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.96%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '64.40'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +7.15%  '
$ws.Range("E42").Value = '  +1.18%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '108.68'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.75%  '
$ws.Range("E44").Value = '  +1.92%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.103'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +2.82%  '
$ws.Range("E46").Value = '  +0.36%  '
$ws.Range("E47").Value = '  +8.49%  '
$ws.Range("E48").Value = '  +1.73%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.20'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +2.25%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.18'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.74%  '
$ws.Range("B51").Value = 'WOONetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.425'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +12.92%  '
